# Update the ref-tracking metric definitions table: several metric rows were
# reworked/renamed/replaced and two new columns' worth of description text
# changed. Re-write the data rows, then re-apply best-fit column widths and
# reset the active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) is unchanged except for the "Optimal" label text reused from
# the shared-string pool - no edit needed there.

# Row 2
$ws.Range("A2").Value = "Trail to Lead"
$ws.Range("B2").Value = "transition_speed"
$ws.Range("C2").Value = "Lead"
$ws.Range("D2").Value = "Transition"
$ws.Range("E2").Value = "Less"
$ws.Range("F2").Value = "Speed at which the lead ref gets to the baseline from the other side of the court at the start of a possession"

# Row 3
$ws.Range("A3").Value = "Time near FT Line Extended"
$ws.Range("B3").Value = "perc_time_near_ft_line_ext"
$ws.Range("C3").Value = "Slot"
$ws.Range("D3").Value = "Halfcourt"
$ws.Range("E3").Value = "More/Undefined"
$ws.Range("F3").Value = "% of time spent within 3 feet of the FT line extended"

# Row 4
$ws.Range("A4").Value = "Trail by 28 Foot Mark"
$ws.Range("B4").Value = "perc_time_by_28_mark"
$ws.Range("C4").Value = "Trail"
$ws.Range("D4").Value = "Halfcourt"
$ws.Range("E4").Value = "More"
$ws.Range("F4").Value = "% of time spent at least 3 feet of the 28 Foot Mark"

# Row 5
$ws.Range("A5").Value = "Trail behind ball - Halfcourt"
$ws.Range("B5").Value = "perc_time_behind_ball_halfcourt"
$ws.Range("C5").Value = "Trail"
$ws.Range("D5").Value = "Halfcourt"
$ws.Range("E5").Value = "More"
$ws.Range("F5").Value = "% of time spent at least 3 feet behind the basketball in the halfcourt"

# Row 6
$ws.Range("A6").Value = "Trail Stay on Play"
$ws.Range("B6").Value = "perc_poss_completed/shifted"
$ws.Range("C6").Value = "Trail"
$ws.Range("D6").Value = "Halfcourt"
$ws.Range("E6").Value = "Less"
$ws.Range("F6").Value = "% of possessions where the trail does not transition until the possession is completed"

# Row 7
$ws.Range("A7").Value = "Trail behind ball - Transition"
$ws.Range("B7").Value = "perc_time_behind_ball_transition"
$ws.Range("C7").Value = "Trail"
$ws.Range("D7").Value = "Transition"
$ws.Range("E7").Value = "More"
$ws.Range("F7").Value = "% of time spent at least 3 feet behind the basketball in transition"

# Row 8
$ws.Range("A8").Value = "Base Lead"
$ws.Range("B8").Value = "perc_time_in_base_position_lead"
$ws.Range("C8").Value = "Lead"
$ws.Range("D8").Value = "Halfcourt"
$ws.Range("E8").Value = "More"
$ws.Range("F8").Value = "% of time spent at least 12 feet out of the rim."

# Re-fit the columns now that longer/shorter text lives in each one: column A
# shrinks, column B grows slightly, and C:F now need best-fit widths of their
# own since they hold real values for the first time.
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 28.498697916666668
$ws.Columns.Item(3).ColumnWidth = 4.166666666666667
$ws.Columns.Item(4).ColumnWidth = 8.498697916666666
$ws.Columns.Item(5).ColumnWidth = 13.998697916666666
$ws.Columns.Item(6).ColumnWidth = 90.99869791666667

# Reset the selection/active cell back to A1 (it had been left on E1).
$ws.Range("A1").Select() | Out-Null
